# Corrected Count attribute for the 'On Campus (excluding Residence Halls)'
# records (rows 2-28) so that the new count = current count - the paired
# 'On Campus (Residence Halls)' record's count (same Sector/Offense/Date,
# 27 rows further down the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 28; $r++) {
    $exclCell = $ws.Cells.Item($r, 5)
    $resCell  = $ws.Cells.Item($r + 27, 5)

    $exclValue = $exclCell.Value2
    $resValue  = $resCell.Value2

    if ($null -eq $exclValue) { $exclValue = 0 }
    if ($null -eq $resValue)  { $resValue  = 0 }

    $exclCell.Value = $exclValue - $resValue
}

# Reflect the final selection / view state recorded in the saved file:
# the window scroll reset to the top and the active cell moved to I7.
$ws.Range("I7").Select()
